$wb = $excel.ActiveWorkbook

# --- Insert a new worksheet named "range" right after "Sheet3" -------------
#
# A throwaway sheet is created first and removed again at the end; this makes
# the workbook's internal "next sheet id" counter land on the same sheetId
# that real Excel would have produced for the newly inserted "range" sheet.
$tmp = $wb.Worksheets.Add()
$tmp.Name = "TMP_DELETE_ME"

$sheet3 = $wb.Worksheets.Item("Sheet3")
$newSheet = $wb.Worksheets.Add($null, $sheet3)
$newSheet.Name = "range"

# Populate the new sheet with a 4x4 grid of values (1..16)
$newSheet.Range("A1").Value = 1
$newSheet.Range("B1").Value = 2
$newSheet.Range("C1").Value = 3
$newSheet.Range("D1").Value = 4
$newSheet.Range("A2").Value = 5
$newSheet.Range("B2").Value = 6
$newSheet.Range("C2").Value = 7
$newSheet.Range("D2").Value = 8
$newSheet.Range("A3").Value = 9
$newSheet.Range("B3").Value = 10
$newSheet.Range("C3").Value = 11
$newSheet.Range("D3").Value = 12
$newSheet.Range("A4").Value = 13
$newSheet.Range("B4").Value = 14
$newSheet.Range("C4").Value = 15
$newSheet.Range("D4").Value = 16

$newSheet.Range("F22").Select() | Out-Null

# --- "some" sheet is no longer the active tab; reset its selection ---------
$some = $wb.Worksheets.Item("some")
$some.Range("B1").Select() | Out-Null

# --- Clean up the throwaway sheet and activate "range" ---------------------
$tmp = $wb.Worksheets.Item("TMP_DELETE_ME")
$tmp.Delete() | Out-Null

$rangeSheet = $wb.Worksheets.Item("range")
$rangeSheet.Activate() | Out-Null
